$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198632001876831
$ws.Range("B1").Value = 2.603037357330322
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.182494401931763
$ws.Range("E1").Value = 1.174451231956482
